# Update values produced by the RandomForest imputation run (terrestrial
# mammals, combination_2_ABCDE/BCE, seed5). The underlying algorithm name
# was refreshed and these cells' results changed accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -13.4257
$ws.Range("E3").Value = 16.1979
$ws.Range("B12").Value = 4.9073
$ws.Range("C14").Value = -13.4879
$ws.Range("C26").Value = -11.59379999999999
$ws.Range("E30").Value = 15.68700000000001
$ws.Range("C31").Value = -13.2147
$ws.Range("B32").Value = 6.534699999999999
$ws.Range("C35").Value = -12.28360000000001
$ws.Range("B36").Value = 9.049700000000007
$ws.Range("C37").Value = -13.2105
$ws.Range("B38").Value = 5.467300000000002
$ws.Range("E44").Value = 16.89829999999999
$ws.Range("C45").Value = -13.85079999999999
$ws.Range("B46").Value = 6.6233
$ws.Range("B54").Value = 4.6351
$ws.Range("B55").Value = 5.310899999999999
$ws.Range("C57").Value = -13.98919999999999
$ws.Range("E58").Value = 16.36980000000001
$ws.Range("B67").Value = 5.832099999999994
$ws.Range("B69").Value = 5.074399999999996
$ws.Range("B72").Value = 5.268900000000002
$ws.Range("E84").Value = 16.45479999999999
$ws.Range("E89").Value = 17.24000000000002
$ws.Range("B91").Value = 5.543199999999998
$ws.Range("E91").Value = 17.90280000000002
$ws.Range("E92").Value = 18.01360000000003
$ws.Range("B99").Value = 4.450299999999999
$ws.Range("C100").Value = -12.438
$ws.Range("C102").Value = -12.2739
$ws.Range("E102").Value = 16.76869999999999
